$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.335.01"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "2.999.73"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'562.49"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'138.81"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("D9").Value = "2.989.78"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +6.97%  "
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "'33.75"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'7.37"
$ws.Range("E16").Value = "  +6.97%  "
$ws.Range("D17").Value = "3.496.01"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.998.03"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "59.291.65"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "'430.89"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").Value = "'13.64"
$ws.Range("E21").Value = "  +2.71%  "
$ws.Range("D22").Value = "'0.723"
$ws.Range("E22").Value = "  +4.94%  "
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").Value = "'80.45"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  +11.27%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").Value = "'25.79"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'0.1000"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  +7.51%  "
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "'49.02"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'8.68"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("D41").Value = "'410.37"
$ws.Range("E41").Value = "  +8.04%  "
$ws.Range("D42").Value = "'0.0354"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").Value = "2.777.12"
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("D47").Value = "'34.81"
$ws.Range("E47").Value = "  +19.79%  "
$ws.Range("D48").Value = "'123.54"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "'2.01"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'23.56"
$ws.Range("E51").Value = "  -0.50%  "
